# Bulk upload fixes: add a "Date Created (Year)*" column right after the
# "filename" column, populate it with the default year value on the
# existing data rows, and leave the selection where the user's cursor
# ended up after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new column before column B; everything from B onward (titles,
# descriptions, documentType, numberOfVolumes, ...) shifts one column right.
$ws.Columns("B").Insert()

# New header + values for the inserted column.
$ws.Range("B1").Value = "Date Created (Year)*"
$ws.Range("B2").Value = 2000
$ws.Range("B3").Value = 2000
$ws.Range("B4").Value = 2000

# Match the page orientation recorded for the sheet after the edit.
$ws.PageSetup.Orientation = 1

# Leave the selection on the last cell touched.
$ws.Range("B4").Select()
